$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79

$ws.Cells.Item($row, 1).Value = "27/06/2024 05:44:41"
$ws.Cells.Item($row, 2).Value = 1
$ws.Cells.Item($row, 3).Value = "NMDC"
$ws.Cells.Item($row, 4).Value = "Nmdc Limited"

# bsecode is stored as text ("526371") rather than a number, unlike the
# other rows in this column, so force text formatting before assignment
# and then drop back to the default style so no stray style index sticks
# to the cell.
$ws.Cells.Item($row, 5).NumberFormat = "@"
$ws.Cells.Item($row, 5).Value = "526371"
$ws.Cells.Item($row, 5).Style = "Normal"

$ws.Cells.Item($row, 6).Value = -0.52
$ws.Cells.Item($row, 7).Value = 247.95
$ws.Cells.Item($row, 8).Value = 5109115
